$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove trailing "*" from the "Prefers sustainable future*" label (A4)
$ws.Range("A4").Value = "Prefers sustainable future"

# Re-run RU 1001 without crop: updated probability values throughout the grid
$ws.Range("B2").Value = 0.693991148793208
$ws.Range("L2").Value = 0.754180914635801
$ws.Range("B3").Value = 0.641501546421187
$ws.Range("L3").Value = 0.759131699735939
$ws.Range("B4").Value = 0.680815013747804
$ws.Range("D4").Value = 0.721220700610068
$ws.Range("E4").Value = 0.703415187125182
$ws.Range("F4").Value = 0.759706650286704
$ws.Range("G4").Value = 0.575101001313061
$ws.Range("H4").Value = 0.731116681894394
$ws.Range("I4").Value = 0.68337742404715
$ws.Range("J4").Value = 0.659704988708263
$ws.Range("K4").Value = 0.759539408850822
$ws.Range("L4").Value = 0.689906588508866
$ws.Range("M4").Value = 0.721853509181885
$ws.Range("N4").Value = 0.618731104894089
$ws.Range("B5").Value = 0.704960018034767
$ws.Range("D5").Value = 0.767261516731427
$ws.Range("E5").Value = 0.757541746745335
$ws.Range("F5").Value = 0.874984008964506
$ws.Range("G5").Value = 0.847206298409435
$ws.Range("H5").Value = 0.842339547896952
$ws.Range("I5").Value = 0.65799487387268
$ws.Range("J5").Value = 0.656332785129309
$ws.Range("K5").Value = 0.703874631903231
$ws.Range("L5").Value = 0.777883926828007
$ws.Range("M5").Value = 0.92816201896394
$ws.Range("N5").Value = 0.562406199574745
$ws.Range("D6").Value = 0.703425024778734
$ws.Range("E6").Value = 0.690418272415403
$ws.Range("F6").Value = 0.817088824971887
$ws.Range("G6").Value = 0.707832398863728
$ws.Range("H6").Value = 0.743080145118413
$ws.Range("I6").Value = 0.682786565268408
$ws.Range("J6").Value = 0.629978960602202
$ws.Range("K6").Value = 0.559364498889858
$ws.Range("L6").Value = 0.559364498889858
$ws.Range("M6").Value = 0.729122438104651
$ws.Range("N6").Value = 0.667115559981111
$ws.Range("D7").Value = 0.716446447911511
$ws.Range("E7").Value = 0.705443804036464
$ws.Range("F7").Value = 0.815943040010246
$ws.Range("G7").Value = 0.639135429260725
$ws.Range("H7").Value = 0.765218616839628
$ws.Range("I7").Value = 0.694592504615592
$ws.Range("J7").Value = 0.566948275489208
$ws.Range("K7").Value = 0.557987381952544
$ws.Range("N7").Value = 0.670153434751546
$ws.Range("D8").Value = 0.438407527923474
$ws.Range("E8").Value = 0.439022530105316
$ws.Range("F8").Value = 0.693097173376912
$ws.Range("H8").Value = 0.510821694970915
$ws.Range("I8").Value = 0.457151305629475
$ws.Range("N8").Value = 0.400676565778138
$ws.Range("B9").Value = 0.591975291077805
$ws.Range("D9").Value = 0.432966587857026
$ws.Range("E9").Value = 0.623667021276442
$ws.Range("F9").Value = 0.765866950290479
$ws.Range("G9").Value = 0.633729019702456
$ws.Range("H9").Value = 0.702512213956321
$ws.Range("I9").Value = 0.575699232924345
$ws.Range("J9").Value = 0.53058178548063
$ws.Range("K9").Value = 0.586217425103406
$ws.Range("L9").Value = 0.572781630922761
$ws.Range("M9").Value = 0.887006592651732
$ws.Range("N9").Value = 0.554098197443099
$ws.Range("B10").Value = 0.0
$ws.Range("C10").Value = 0.0
$ws.Range("D10").Value = 0.0
$ws.Range("E10").Value = 0.0
$ws.Range("F10").Value = 0.0
$ws.Range("G10").Value = 0.0
$ws.Range("H10").Value = 0.0
$ws.Range("I10").Value = 0.0
$ws.Range("J10").Value = 0.0
$ws.Range("K10").Value = 0.0
$ws.Range("L10").Value = 0.0
$ws.Range("M10").Value = 0.0
$ws.Range("N10").Value = 0.0
